$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 383, pushing existing rows 383.. down to 384..
$ws.Rows.Item(383).Insert()

# Populate the newly inserted row 383 with the new record's data
$ws.Cells.Item(383, 1).Value = 5
$ws.Cells.Item(383, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(383, 3).Value = "Maule"
$ws.Cells.Item(383, 4).Value = 45211
$ws.Cells.Item(383, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(383, 5).Value = 7
$ws.Cells.Item(383, 6).Value = "Fruta"
$ws.Cells.Item(383, 7).Value = 100108
$ws.Cells.Item(383, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(383, 9).Value = 100108005
$ws.Cells.Item(383, 10).Value = "Piña"
$ws.Cells.Item(383, 11).Value = "Caramelo"
$ws.Cells.Item(383, 12).Value = "Segunda"
$ws.Cells.Item(383, 13).Value = 180
$ws.Cells.Item(383, 14).Value = 22000
$ws.Cells.Item(383, 15).Value = 22000
$ws.Cells.Item(383, 16).Value = 22000
$ws.Cells.Item(383, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(383, 18).Value = "Ecuador"
$ws.Cells.Item(383, 19).Value = 1571
$ws.Cells.Item(383, 20).Value = 14
